{"js": "// Locate the existing list paragraph that ends the list (\"Sublime text\")\n// and add a new list item \"Adobe Flash\" right after it, matching the\n// list's paragraph style / numbering / language formatting.\nconst body = context.document.body;\n\nconst results = body.search(\"Sublime text\", { matchCase: false, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find paragraph containing \"Sublime text\"');\n}\n\nconst anchorRange = results.items[0];\nconst anchorParagraph = anchorRange.paragraphs.getFirst();\n\n// Insert a new paragraph right after the anchor paragraph; insertParagraph\n// on a paragraph carries over that paragraph's style/list formatting.\nconst newParagraph = anchorParagraph.insertParagraph(\"Adobe Flash\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new list item \"Adobe Flash\" right after the existing \"Sublime text\"\n# list item, matching its paragraph style / numbering / language formatting.\n$d = $word.ActiveDocument\n\n# Find the (list) paragraph that holds \"Sublime text\" and remember its index.\n$targetIndex = -1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Sublime text*\") {\n        $targetIndex = $p.Index\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find paragraph containing 'Sublime text'\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n\n# InsertParagraphAfter creates a new paragraph that inherits the formatting\n# (style/list numbering) of $target.\n$target.Range.InsertParagraphAfter()\n\n# Re-fetch the freshly created paragraph (immediately after $target) and set\n# its text.\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Adobe Flash\"\n\n$d.Save()\n"}
